$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 10428.353  # H51: 10225.823 -> 10428.353
$ws.Cells.Item(51, 9).Value = 22332.666  # I51: 17388.25 -> 22332.666
$ws.Cells.Item(51, 10).Value = 7877.4287  # J51: 8022 -> 7877.4287
$ws.Cells.Item(51, 11).Value = 22332.666  # K51: 17388.25 -> 22332.666
$ws.Cells.Item(51, 12).Value = 7877.4287  # L51: 8022 -> 7877.4287
$ws.Cells.Item(51, 13).Value = -21848.666  # M51: -16904.25 -> -21848.666
$ws.Cells.Item(51, 14).Value = -8845.4287  # N51: -8990 -> -8845.4287

$ws.Cells.Item(80, 8).Value = 23198.363  # H80: 24307.38 -> 23198.363
$ws.Cells.Item(80, 9).Value = 7565.857  # I80: 7583.7144 -> 7565.857
$ws.Cells.Item(80, 10).Value = 50555.25  # J80: 57754.715 -> 50555.25
$ws.Cells.Item(80, 11).Value = 22697.571  # K80: 22751.1432 -> 22697.571
$ws.Cells.Item(80, 12).Value = 151665.75  # L80: 173264.145 -> 151665.75
$ws.Cells.Item(80, 13).Value = -21699.571  # M80: -21753.1432 -> -21699.571
$ws.Cells.Item(80, 14).Value = -153661.75  # N80: -175260.145 -> -153661.75

$ws.Cells.Item(83, 8).Value = 23198.363  # H83: 24307.38 -> 23198.363
$ws.Cells.Item(83, 9).Value = 7565.857  # I83: 7583.7144 -> 7565.857
$ws.Cells.Item(83, 10).Value = 50555.25  # J83: 57754.715 -> 50555.25
$ws.Cells.Item(83, 11).Value = 68092.713  # K83: 68253.4296 -> 68092.713
$ws.Cells.Item(83, 12).Value = 454997.25  # L83: 519792.4349999999 -> 454997.25
$ws.Cells.Item(83, 13).Value = -63100.713  # M83: -63261.4296 -> -63100.713
$ws.Cells.Item(83, 14).Value = -464981.25  # N83: -529776.4349999999 -> -464981.25

$ws.Cells.Item(113, 8).Value = 35953708  # H113: 33956390 -> 35953708
$ws.Cells.Item(113, 9).Value = 22224406  # I113: 18520672 -> 22224406
$ws.Cells.Item(113, 11).Value = 22224406  # K113: 18520672 -> 22224406
$ws.Cells.Item(113, 13).Value = -22221152  # M113: -18517418 -> -22221152

$ws.Cells.Item(137, 8).Value = 3632.2  # H137: 2447.3914 -> 3632.2
$ws.Cells.Item(137, 9).Value = 5364.3335  # I137: 2524.4285 -> 5364.3335
$ws.Cells.Item(137, 10).Value = 2477.4443  # J137: 2327.5557 -> 2477.4443
$ws.Cells.Item(137, 11).Value = 16093.0005  # K137: 7573.2855 -> 16093.0005
$ws.Cells.Item(137, 12).Value = 7432.3329  # L137: 6982.6671 -> 7432.3329
$ws.Cells.Item(137, 13).Value = -13543.0005  # M137: -5023.2855 -> -13543.0005
$ws.Cells.Item(137, 14).Value = -12532.3329  # N137: -12082.6671 -> -12532.3329

$ws.Cells.Item(138, 8).Value = 3282.8928  # H138: 3288.3103 -> 3282.8928
$ws.Cells.Item(138, 9).Value = 894.88635  # I138: 922.2857 -> 894.88635
$ws.Cells.Item(138, 10).Value = 5909.7  # J138: 5496.6 -> 5909.7
$ws.Cells.Item(138, 11).Value = 2684.65905  # K138: 2766.8571 -> 2684.65905
$ws.Cells.Item(138, 12).Value = 17729.1  # L138: 16489.8 -> 17729.1
$ws.Cells.Item(138, 13).Value = 2455.34095  # M138: 2373.1429 -> 2455.34095
$ws.Cells.Item(138, 14).Value = -28009.1  # N138: -26769.8 -> -28009.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 3940.2144  # H2: 3708.8 -> 3940.2144
$ws.Cells.Item(2, 9).Value = 1132  # I2: 1058.3334 -> 1132
$ws.Cells.Item(2, 11).Value = 1132  # K2: 1058.3334 -> 1132
$ws.Cells.Item(2, 13).Value = -1019  # M2: -945.3334 -> -1019

$ws.Cells.Item(32, 8).Value = 1440521.6  # H32: 1457270.9 -> 1440521.6
$ws.Cells.Item(32, 9).Value = 1509093.6  # I32: 1509093.8 -> 1509093.6
$ws.Cells.Item(32, 10).Value = 17653.25  # J32: 23504.334 -> 17653.25
$ws.Cells.Item(32, 11).Value = 1509093.6  # K32: 1509093.8 -> 1509093.6
$ws.Cells.Item(32, 12).Value = 17653.25  # L32: 23504.334 -> 17653.25
$ws.Cells.Item(32, 13).Value = -1508806.6  # M32: -1508806.8 -> -1508806.6
$ws.Cells.Item(32, 14).Value = -18227.25  # N32: -24078.334 -> -18227.25

$ws.Cells.Item(45, 8).Value = 11465.625  # H45: 11465.875 -> 11465.625
$ws.Cells.Item(45, 10).Value = 15865.4  # J45: 15865.8 -> 15865.4
$ws.Cells.Item(45, 12).Value = 15865.4  # L45: 15865.8 -> 15865.4
$ws.Cells.Item(45, 14).Value = -16619.4  # N45: -16619.8 -> -16619.4

$ws.Cells.Item(116, 8).Value = 3940.2144  # H116: 3708.8 -> 3940.2144
$ws.Cells.Item(116, 9).Value = 1132  # I116: 1058.3334 -> 1132
$ws.Cells.Item(116, 11).Value = 1132  # K116: 1058.3334 -> 1132
$ws.Cells.Item(116, 13).Value = 1162  # M116: 1235.6666 -> 1162

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 3940.2144  # H3: 3708.8 -> 3940.2144
$ws.Cells.Item(3, 9).Value = 1132  # I3: 1058.3334 -> 1132
$ws.Cells.Item(3, 11).Value = 1132  # K3: 1058.3334 -> 1132
$ws.Cells.Item(3, 13).Value = -1018  # M3: -944.3334 -> -1018

$ws.Cells.Item(24, 8).Value = 5500  # H24: 3972 -> 5500
$ws.Cells.Item(24, 9).Value = 1000  # I24: 958 -> 1000
$ws.Cells.Item(24, 11).Value = 1000  # K24: 958 -> 1000
$ws.Cells.Item(24, 13).Value = -765  # M24: -723 -> -765

$ws.Cells.Item(80, 10).Value = 303.33334  # J80: 306.66666 -> 303.33334
$ws.Cells.Item(80, 12).Value = 303.33334  # L80: 306.66666 -> 303.33334
$ws.Cells.Item(80, 14).Value = -2299.33334  # N80: -2302.66666 -> -2299.33334

$ws.Cells.Item(83, 10).Value = 303.33334  # J83: 306.66666 -> 303.33334
$ws.Cells.Item(83, 12).Value = 1516.6667  # L83: 1533.3333 -> 1516.6667
$ws.Cells.Item(83, 14).Value = -11500.6667  # N83: -11517.3333 -> -11500.6667

$ws.Cells.Item(88, 8).Value = 35360  # H88: 0 -> 35360
$ws.Cells.Item(88, 10).Value = 35360  # J88: 0 -> 35360
$ws.Cells.Item(88, 12).Value = 35360  # L88: 0 -> 35360
$ws.Cells.Item(88, 14).Value = -36172  # N88: None -> -36172

$ws.Cells.Item(91, 8).Value = 35360  # H91: 0 -> 35360
$ws.Cells.Item(91, 10).Value = 35360  # J91: 0 -> 35360
$ws.Cells.Item(91, 12).Value = 35360  # L91: 0 -> 35360
$ws.Cells.Item(91, 14).Value = -38168  # N91: None -> -38168

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6624.306  # H31: 7008.298 -> 6624.306
$ws.Cells.Item(31, 9).Value = 2263.5  # I31: 2536.3845 -> 2263.5
$ws.Cells.Item(31, 10).Value = 11553.913  # J31: 12544.952 -> 11553.913
$ws.Cells.Item(31, 11).Value = 2263.5  # K31: 2536.3845 -> 2263.5
$ws.Cells.Item(31, 12).Value = 11553.913  # L31: 12544.952 -> 11553.913
$ws.Cells.Item(31, 13).Value = -1968.5  # M31: -2241.3845 -> -1968.5
$ws.Cells.Item(31, 14).Value = -12143.913  # N31: -13134.952 -> -12143.913

$ws.Cells.Item(34, 8).Value = 6624.306  # H34: 7008.298 -> 6624.306
$ws.Cells.Item(34, 9).Value = 2263.5  # I34: 2536.3845 -> 2263.5
$ws.Cells.Item(34, 10).Value = 11553.913  # J34: 12544.952 -> 11553.913
$ws.Cells.Item(34, 11).Value = 2263.5  # K34: 2536.3845 -> 2263.5
$ws.Cells.Item(34, 12).Value = 11553.913  # L34: 12544.952 -> 11553.913
$ws.Cells.Item(34, 13).Value = -2061.5  # M34: -2334.3845 -> -2061.5
$ws.Cells.Item(34, 14).Value = -11957.913  # N34: -12948.952 -> -11957.913

$ws.Cells.Item(99, 8).Value = 13174.875  # H99: 13724.75 -> 13174.875
$ws.Cells.Item(99, 9).Value = 19133  # I99: 26899.5 -> 19133
$ws.Cells.Item(99, 10).Value = 9600  # J99: 9333.166999999999 -> 9600
$ws.Cells.Item(99, 11).Value = 19133  # K99: 26899.5 -> 19133
$ws.Cells.Item(99, 12).Value = 9600  # L99: 9333.166999999999 -> 9600
$ws.Cells.Item(99, 13).Value = -17635  # M99: -25401.5 -> -17635
$ws.Cells.Item(99, 14).Value = -12596  # N99: -12329.167 -> -12596

$ws.Cells.Item(126, 8).Value = 13174.875  # H126: 13724.75 -> 13174.875
$ws.Cells.Item(126, 9).Value = 19133  # I126: 26899.5 -> 19133
$ws.Cells.Item(126, 10).Value = 9600  # J126: 9333.166999999999 -> 9600
$ws.Cells.Item(126, 11).Value = 57399  # K126: 80698.5 -> 57399
$ws.Cells.Item(126, 12).Value = 28800  # L126: 27999.501 -> 28800
$ws.Cells.Item(126, 13).Value = -54929  # M126: -78228.5 -> -54929
$ws.Cells.Item(126, 14).Value = -33740  # N126: -32939.501 -> -33740

$ws.Cells.Item(134, 8).Value = 7030.161  # H134: 7187.4136 -> 7030.161
$ws.Cells.Item(134, 9).Value = 1751.8  # I134: 1835.3334 -> 1751.8
$ws.Cells.Item(134, 10).Value = 9543.666999999999  # J134: 9595.85 -> 9543.666999999999
$ws.Cells.Item(134, 11).Value = 5255.4  # K134: 5506.0002 -> 5255.4
$ws.Cells.Item(134, 12).Value = 28631.001  # L134: 28787.55 -> 28631.001
$ws.Cells.Item(134, 13).Value = -2720.4  # M134: -2971.0002 -> -2720.4
$ws.Cells.Item(134, 14).Value = -33701.001  # N134: -33857.55 -> -33701.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 15151827  # H14: 13889196 -> 15151827
$ws.Cells.Item(14, 9).Value = 15151827  # I14: 13889196 -> 15151827
$ws.Cells.Item(14, 11).Value = 45455481  # K14: 41667588 -> 45455481
$ws.Cells.Item(14, 13).Value = -45455308  # M14: -41667415 -> -45455308

$ws.Cells.Item(68, 8).Value = 5956.6  # H68: 5954.6 -> 5956.6
$ws.Cells.Item(68, 9).Value = 1865  # I68: 1862.5 -> 1865
$ws.Cells.Item(68, 11).Value = 5595  # K68: 5587.5 -> 5595
$ws.Cells.Item(68, 13).Value = -4784  # M68: -4776.5 -> -4784

$ws.Cells.Item(71, 8).Value = 5956.6  # H71: 5954.6 -> 5956.6
$ws.Cells.Item(71, 9).Value = 1865  # I71: 1862.5 -> 1865
$ws.Cells.Item(71, 11).Value = 16785  # K71: 16762.5 -> 16785
$ws.Cells.Item(71, 13).Value = -12729  # M71: -12706.5 -> -12729

$ws.Cells.Item(98, 8).Value = 3569.5833  # H98: 3521.3845 -> 3569.5833
$ws.Cells.Item(98, 10).Value = 3856.625  # J98: 3755.111 -> 3856.625
$ws.Cells.Item(98, 12).Value = 11569.875  # L98: 11265.333 -> 11569.875
$ws.Cells.Item(98, 14).Value = -14565.875  # N98: -14261.333 -> -14565.875

$ws.Cells.Item(129, 8).Value = 15152920  # H129: 90910500 -> 15152920
$ws.Cells.Item(129, 10).Value = 23811238  # J129: 142858850 -> 23811238
$ws.Cells.Item(129, 12).Value = 71433714  # L129: 428576550 -> 71433714
$ws.Cells.Item(129, 14).Value = -71443714  # N129: -428586550 -> -71443714

$ws.Cells.Item(131, 8).Value = 1333.9445  # H131: 1371.2941 -> 1333.9445
$ws.Cells.Item(131, 9).Value = 906.1429000000001  # I131: 940.6667 -> 906.1429000000001
$ws.Cells.Item(131, 11).Value = 2718.4287  # K131: 2822.0001 -> 2718.4287
$ws.Cells.Item(131, 13).Value = 2321.5713  # M131: 2217.9999 -> 2321.5713

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 5053.5454  # H80: 5053.636 -> 5053.5454
$ws.Cells.Item(80, 9).Value = 6024.25  # I80: 7033.3335 -> 6024.25
$ws.Cells.Item(80, 10).Value = 4498.857  # J80: 4311.25 -> 4498.857
$ws.Cells.Item(80, 11).Value = 6024.25  # K80: 7033.3335 -> 6024.25
$ws.Cells.Item(80, 12).Value = 4498.857  # L80: 4311.25 -> 4498.857
$ws.Cells.Item(80, 13).Value = -5026.25  # M80: -6035.3335 -> -5026.25
$ws.Cells.Item(80, 14).Value = -6494.857  # N80: -6307.25 -> -6494.857

$ws.Cells.Item(83, 8).Value = 5053.5454  # H83: 5053.636 -> 5053.5454
$ws.Cells.Item(83, 9).Value = 6024.25  # I83: 7033.3335 -> 6024.25
$ws.Cells.Item(83, 10).Value = 4498.857  # J83: 4311.25 -> 4498.857
$ws.Cells.Item(83, 11).Value = 30121.25  # K83: 35166.6675 -> 30121.25
$ws.Cells.Item(83, 12).Value = 22494.285  # L83: 21556.25 -> 22494.285
$ws.Cells.Item(83, 13).Value = -25129.25  # M83: -30174.6675 -> -25129.25
$ws.Cells.Item(83, 14).Value = -32478.285  # N83: -31540.25 -> -32478.285

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6421.5  # H7: 6034.4165 -> 6421.5
$ws.Cells.Item(7, 9).Value = 4675.4546  # I7: 4652.1665 -> 4675.4546
$ws.Cells.Item(7, 10).Value = 8555.556  # J7: 7416.6665 -> 8555.556
$ws.Cells.Item(7, 11).Value = 4675.4546  # K7: 4652.1665 -> 4675.4546
$ws.Cells.Item(7, 12).Value = 8555.556  # L7: 7416.6665 -> 8555.556
$ws.Cells.Item(7, 13).Value = -4563.4546  # M7: -4540.1665 -> -4563.4546
$ws.Cells.Item(7, 14).Value = -8779.556  # N7: -7640.6665 -> -8779.556

$ws.Cells.Item(40, 8).Value = 8141.9165  # H40: 7986.7334 -> 8141.9165
$ws.Cells.Item(40, 9).Value = 6939.6  # I40: 7159.8 -> 6939.6
$ws.Cells.Item(40, 10).Value = 9000.714  # J40: 8400.200000000001 -> 9000.714
$ws.Cells.Item(40, 11).Value = 6939.6  # K40: 7159.8 -> 6939.6
$ws.Cells.Item(40, 12).Value = 9000.714  # L40: 8400.200000000001 -> 9000.714
$ws.Cells.Item(40, 13).Value = -6803.6  # M40: -7023.8 -> -6803.6
$ws.Cells.Item(40, 14).Value = -9272.714  # N40: -8672.200000000001 -> -9272.714

$ws.Cells.Item(93, 8).Value = 9089  # H93: 6986 -> 9089
$ws.Cells.Item(93, 9).Value = 8450.5  # I93: 3488.6667 -> 8450.5
$ws.Cells.Item(93, 10).Value = 9271.429  # J93: 10483.333 -> 9271.429
$ws.Cells.Item(93, 11).Value = 8450.5  # K93: 3488.6667 -> 8450.5
$ws.Cells.Item(93, 12).Value = 9271.429  # L93: 10483.333 -> 9271.429
$ws.Cells.Item(93, 13).Value = -7202.5  # M93: -2240.6667 -> -7202.5
$ws.Cells.Item(93, 14).Value = -11767.429  # N93: -12979.333 -> -11767.429

$ws.Cells.Item(100, 8).Value = 4469.5  # H100: 4582.6924 -> 4469.5
$ws.Cells.Item(100, 9).Value = 3457.2  # I100: 3508.2222 -> 3457.2
$ws.Cells.Item(100, 11).Value = 3457.2  # K100: 3508.2222 -> 3457.2
$ws.Cells.Item(100, 13).Value = -2916.2  # M100: -2967.2222 -> -2916.2

$ws.Cells.Item(122, 8).Value = 4251.1665  # H122: 4273.736 -> 4251.1665
$ws.Cells.Item(122, 9).Value = 3464.9756  # I122: 3475.225 -> 3464.9756
$ws.Cells.Item(122, 11).Value = 10394.9268  # K122: 10425.675 -> 10394.9268
$ws.Cells.Item(122, 13).Value = -7944.926800000001  # M122: -7975.674999999999 -> -7944.926800000001

$ws.Cells.Item(126, 8).Value = 6421.5  # H126: 6034.4165 -> 6421.5
$ws.Cells.Item(126, 9).Value = 4675.4546  # I126: 4652.1665 -> 4675.4546
$ws.Cells.Item(126, 10).Value = 8555.556  # J126: 7416.6665 -> 8555.556
$ws.Cells.Item(126, 11).Value = 14026.3638  # K126: 13956.4995 -> 14026.3638
$ws.Cells.Item(126, 12).Value = 25666.668  # L126: 22249.9995 -> 25666.668
$ws.Cells.Item(126, 13).Value = -11556.3638  # M126: -11486.4995 -> -11556.3638
$ws.Cells.Item(126, 14).Value = -30606.668  # N126: -27189.9995 -> -30606.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 17500  # H54: 16900 -> 17500
$ws.Cells.Item(54, 10).Value = 20000  # J54: 19750 -> 20000
$ws.Cells.Item(54, 12).Value = 20000  # L54: 19750 -> 20000
$ws.Cells.Item(54, 14).Value = -21040  # N54: -20790 -> -21040

$ws.Cells.Item(56, 8).Value = 0  # H56: 33330 -> 0
$ws.Cells.Item(56, 10).Value = 0  # J56: 33330 -> 0
$ws.Cells.Item(56, 12).Value = 0  # L56: 33330 -> 0
$ws.Cells.Item(56, 14).ClearContents()  # N56: remove (was -34758)

$ws.Cells.Item(96, 8).Value = 2461.5  # H96: 2798.1667 -> 2461.5
$ws.Cells.Item(96, 9).Value = 2527.4285  # I96: 2857.8 -> 2527.4285
$ws.Cells.Item(96, 10).Value = 2000  # J96: 2500 -> 2000
$ws.Cells.Item(96, 11).Value = 2527.4285  # K96: 2857.8 -> 2527.4285
$ws.Cells.Item(96, 12).Value = 2000  # L96: 2500 -> 2000
$ws.Cells.Item(96, 13).Value = -1154.4285  # M96: -1484.8 -> -1154.4285
$ws.Cells.Item(96, 14).Value = -4746  # N96: -5246 -> -4746

$ws.Cells.Item(113, 8).Value = 678.25  # H113: 826.58826 -> 678.25
$ws.Cells.Item(113, 9).Value = 565.8182  # I113: 785.3333 -> 565.8182
$ws.Cells.Item(113, 11).Value = 1697.4546  # K113: 2355.9999 -> 1697.4546
$ws.Cells.Item(113, 13).Value = 472.5454  # M113: -185.9998999999998 -> 472.5454
